# Add DataTables (DT) internationalization strings to the translations sheet.
# This inserts 8 new rows right before the "Database columns and table names"
# section (previously row 47), pushing that section and the rows after it
# down by 8 rows, then fills the newly inserted rows with the new id /
# description / English / French values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 blank rows at row 47 (keeps formatting of the row that was there,
# which matches the target diff where the inserted rows are plain/unstyled).
$ws.Range("A47:A54").EntireRow.Insert()

$ws.Range("A47").Value = "tbl_info"
$ws.Range("B47").Value = "DT internationalization"
$ws.Range("C47").Value = " _TOTAL_ records total"
$ws.Range("D47").Value = "_TOTAL_ éléments au total"

$ws.Range("A48").Value = "tbl_info_empty"
$ws.Range("B48").Value = "DT internationalization"
$ws.Range("C48").Value = "No records to show!"
$ws.Range("D48").Value = "Rien à voir ici!"

$ws.Range("A49").Value = "tbl_prev"
$ws.Range("B49").Value = "DT internationalization"
$ws.Range("C49").Value = "Previous"
$ws.Range("D49").Value = "Précédent"

$ws.Range("A50").Value = "tbl_next"
$ws.Range("B50").Value = "DT internationalization"
$ws.Range("C50").Value = "Next"
$ws.Range("D50").Value = "Prochain"

$ws.Range("A51").Value = "tbl_search"
$ws.Range("B51").Value = "DT internationalization"
$ws.Range("C51").Value = "Search:"
$ws.Range("D51").Value = "Recherche"

$ws.Range("A52").Value = "tbl_length"
$ws.Range("B52").Value = "DT internationalization"
$ws.Range("C52").Value = "Show _MENU_ records"
$ws.Range("D52").Value = "Montrer _MENU_ éléments"

$ws.Range("A53").Value = "tbl_filtered"
$ws.Range("B53").Value = "DT internationalization"
$ws.Range("C53").Value = "(filtered from _MAX_ records)"
$ws.Range("D53").Value = "(filtrées de _MAX_ éléments)"

$ws.Range("A54").Value = "tbl_zero"
$ws.Range("B54").Value = "DT internationalization"
$ws.Range("C54").Value = "No matching records found"
$ws.Range("D54").Value = "Aucunes données existantes"

# Keep the tab selected / update the active cell to mirror the sheetView
# selection recorded after the edit (row 48, column D) and bring the newly
# edited area into view.
$ws.Range("A41").Select()
$ws.Range("D48").Select()
